$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (rows 3-75) ---
$ws.Range("G3").Value = 0.0113848495233494
$ws.Range("G4").Value = 0.0113848495233494
$ws.Range("G5").Value = 0.0047000715365022
$ws.Range("G6").Value = 0.0047000715365022
$ws.Range("F7").Value = 0.001
$ws.Range("G7").Value = 0.0026133997844927
$ws.Range("F8").Value = 0.001
$ws.Range("G8").Value = 0.0026133997844927
$ws.Range("F9").Value = 0.00525
$ws.Range("G9").Value = 0.0113770347015022
$ws.Range("H9").Value = 0.0702
$ws.Range("I9").Value = 0.05124
$ws.Range("L9").Value = 0.00335
$ws.Range("M9").Value = 0.0175
$ws.Range("F10").Value = 0.00525
$ws.Range("G10").Value = 0.0113770347015022
$ws.Range("H10").Value = 0.0702
$ws.Range("I10").Value = 0.05124
$ws.Range("L10").Value = 0.00335
$ws.Range("M10").Value = 0.0175
$ws.Range("G12").Value = 0.0128346823336127
$ws.Range("G13").Value = 0.0128346823336127
$ws.Range("G14").Value = 0.0043243214908769
$ws.Range("G15").Value = 0.0043243214908769
$ws.Range("G16").Value = 0.0020987659679818
$ws.Range("L16").Value = 0.00046
$ws.Range("G17").Value = 0.0020987659679818
$ws.Range("L17").Value = 0.00046
$ws.Range("F18").Value = 0.0047
$ws.Range("G18").Value = 0.0106370152232541
$ws.Range("H18").Value = 0.0702
$ws.Range("I18").Value = 0.04956
$ws.Range("L18").Value = 0.00366
$ws.Range("F19").Value = 0.0047
$ws.Range("G19").Value = 0.0106370152232541
$ws.Range("H19").Value = 0.0702
$ws.Range("I19").Value = 0.04956
$ws.Range("L19").Value = 0.00366
$ws.Range("G21").Value = 0.0147401837743907
$ws.Range("G22").Value = 0.0147401837743907
$ws.Range("F23").Value = 0.001
$ws.Range("G23").Value = 0.0044111112851871
$ws.Range("L23").Value = 0.00074
$ws.Range("F24").Value = 0.001
$ws.Range("G24").Value = 0.0044111112851871
$ws.Range("L24").Value = 0.00074
$ws.Range("F25").Value = 0.001
$ws.Range("G25").Value = 0.0021641316750224
$ws.Range("L25").Value = 0.00055
$ws.Range("F26").Value = 0.001
$ws.Range("G26").Value = 0.0021641316750224
$ws.Range("L26").Value = 0.00055
$ws.Range("F27").Value = 0.00475
$ws.Range("G27").Value = 0.009710239926219
$ws.Range("H27").Value = 0.0702
$ws.Range("I27").Value = 0.04326
$ws.Range("L27").Value = 0.00386
$ws.Range("F28").Value = 0.00475
$ws.Range("G28").Value = 0.009710239926219
$ws.Range("H28").Value = 0.0702
$ws.Range("I28").Value = 0.04326
$ws.Range("L28").Value = 0.00386
$ws.Range("F32").Value = 0.00163
$ws.Range("G32").Value = 0.436783820531418
$ws.Range("L32").Value = 0.00154
$ws.Range("M32").Value = 0.01123
$ws.Range("F33").Value = 0.00163
$ws.Range("G33").Value = 0.436783820531418
$ws.Range("L33").Value = 0.00154
$ws.Range("M33").Value = 0.01123
$ws.Range("F34").Value = 0.00073
$ws.Range("G34").Value = 0.0013090121203059
$ws.Range("L34").Value = 0.00071
$ws.Range("F35").Value = 0.00073
$ws.Range("G35").Value = 0.0013090121203059
$ws.Range("L35").Value = 0.00071
$ws.Range("F36").Value = 0.00584
$ws.Range("G36").Value = 1.13406246598869
$ws.Range("I36").Value = 0.0522
$ws.Range("L36").Value = 0.0048
$ws.Range("M36").Value = 0.01185
$ws.Range("F37").Value = 0.00584
$ws.Range("G37").Value = 1.13406246598869
$ws.Range("I37").Value = 0.0522
$ws.Range("L37").Value = 0.0048
$ws.Range("M37").Value = 0.01185
$ws.Range("F41").Value = 0.0022
$ws.Range("G41").Value = 0.426576543357142
$ws.Range("I41").Value = 0.02518
$ws.Range("L41").Value = 0.00165
$ws.Range("M41").Value = 0.01179
$ws.Range("N41").Value = 0.01833
$ws.Range("F42").Value = 0.0022
$ws.Range("G42").Value = 0.426576543357142
$ws.Range("I42").Value = 0.02518
$ws.Range("L42").Value = 0.00165
$ws.Range("M42").Value = 0.01179
$ws.Range("N42").Value = 0.01833
$ws.Range("F43").Value = 0.00076
$ws.Range("G43").Value = 0.0012102350402381
$ws.Range("L43").Value = 0.00072
$ws.Range("M43").Value = 0.0019
$ws.Range("F44").Value = 0.00076
$ws.Range("G44").Value = 0.0012102350402381
$ws.Range("L44").Value = 0.00072
$ws.Range("M44").Value = 0.0019
$ws.Range("G45").Value = 1.10886824001553
$ws.Range("M45").Value = 0.01408
$ws.Range("N45").Value = 0.01993
$ws.Range("G46").Value = 1.10886824001553
$ws.Range("M46").Value = 0.01408
$ws.Range("N46").Value = 0.01993
$ws.Range("F50").Value = 0.00325
$ws.Range("G50").Value = 0.43735529010491
$ws.Range("I50").Value = 0.02717
$ws.Range("M50").Value = 0.01186
$ws.Range("N50").Value = 0.01942
$ws.Range("F51").Value = 0.00325
$ws.Range("G51").Value = 0.43735529010491
$ws.Range("I51").Value = 0.02717
$ws.Range("M51").Value = 0.01186
$ws.Range("N51").Value = 0.01942
$ws.Range("G52").Value = 0.0015247305300927
$ws.Range("L52").Value = 0.0008899999999999999
$ws.Range("G53").Value = 0.0015247305300927
$ws.Range("L53").Value = 0.0008899999999999999
$ws.Range("G54").Value = 1.0639713565819
$ws.Range("M54").Value = 0.01363
$ws.Range("N54").Value = 0.02342
$ws.Range("G55").Value = 1.0639713565819
$ws.Range("M55").Value = 0.01363
$ws.Range("N55").Value = 0.02342
$ws.Range("F59").Value = 0.00372
$ws.Range("G59").Value = 0.427653725671211
$ws.Range("I59").Value = 0.02702
$ws.Range("L59").Value = 0.00568
$ws.Range("M59").Value = 0.01248
$ws.Range("F60").Value = 0.00372
$ws.Range("G60").Value = 0.427653725671211
$ws.Range("I60").Value = 0.02702
$ws.Range("L60").Value = 0.00568
$ws.Range("M60").Value = 0.01248
$ws.Range("F61").Value = 0.00103
$ws.Range("G61").Value = 0.0015219531658587
$ws.Range("F62").Value = 0.00103
$ws.Range("G62").Value = 0.0015219531658587
$ws.Range("G63").Value = 1.06422070458161
$ws.Range("L63").Value = 0.009169999999999999
$ws.Range("M63").Value = 0.01673
$ws.Range("G64").Value = 1.06422070458161
$ws.Range("L64").Value = 0.009169999999999999
$ws.Range("M64").Value = 0.01673
$ws.Range("F68").Value = 0.00361
$ws.Range("G68").Value = 0.38305399641797
$ws.Range("I68").Value = 0.0299
$ws.Range("L68").Value = 0.00568
$ws.Range("M68").Value = 0.01266
$ws.Range("F69").Value = 0.00361
$ws.Range("G69").Value = 0.38305399641797
$ws.Range("I69").Value = 0.0299
$ws.Range("L69").Value = 0.00568
$ws.Range("M69").Value = 0.01266
$ws.Range("F70").Value = 0.00118
$ws.Range("G70").Value = 0.0014798114806265
$ws.Range("L70").Value = 0.0011
$ws.Range("M70").Value = 0.00192
$ws.Range("F71").Value = 0.00118
$ws.Range("G71").Value = 0.0014798114806265
$ws.Range("L71").Value = 0.0011
$ws.Range("M71").Value = 0.00192
$ws.Range("F72").Value = 0.00831
$ws.Range("G72").Value = 0.982507565213656
$ws.Range("M72").Value = 0.01555
$ws.Range("N72").Value = 0.02092
$ws.Range("F73").Value = 0.00831
$ws.Range("G73").Value = 0.982507565213656
$ws.Range("M73").Value = 0.01555
$ws.Range("N73").Value = 0.02092
$ws.Range("G74").Value = 0.0523275496330669
$ws.Range("G75").Value = 0.0523275496330669

# --- Add new rows 78-90 for 2019 - 2023 period ---
$ws.Range("A78").Value = "Moawhango at Waiouru"
$ws.Range("B78").Value = "Chlorophyll A (92nd Percentile)"
$ws.Range("C78").Value = "D"
$ws.Range("D78").Value = "2019 - 2023"
$ws.Range("E78").Value = "RepSite"
$ws.Range("F78").Value = 157.5
$ws.Range("G78").Value = 168.663793103448
$ws.Range("H78").Value = 440
$ws.Range("I78").Value = 308
$ws.Range("L78").Value = 180
$ws.Range("M78").Value = 233.2
$ws.Range("N78").Value = 275
$ws.Range("O78").Value = 1838908
$ws.Range("P78").Value = 5631468
$ws.Range("Q78").Value = "Ruapehu District"
$ws.Range("R78").Value = "Rangitīkei-Turakina"
$ws.Range("S78").Value = "Middle Rangitikei"
$ws.Range("T78").Value = "Rang_2d"
$ws.Range("U78").Value = "mg chl-a /m2"
$ws.Range("A79").Value = "Moawhango at Waiouru"
$ws.Range("B79").Value = "DRP (95th Percentile)"
$ws.Range("C79").Value = "B"
$ws.Range("D79").Value = "2019 - 2023"
$ws.Range("E79").Value = "RepSite"
$ws.Range("F79").Value = 0.008999999999999999
$ws.Range("G79").Value = 0.0125192307692308
$ws.Range("H79").Value = 0.078
$ws.Range("I79").Value = 0.0259
$ws.Range("L79").Value = 0.008500000000000001
$ws.Range("M79").Value = 0.018
$ws.Range("N79").Value = 0.02368
$ws.Range("O79").Value = 1838908
$ws.Range("P79").Value = 5631468
$ws.Range("Q79").Value = "Ruapehu District"
$ws.Range("R79").Value = "Rangitīkei-Turakina"
$ws.Range("S79").Value = "Middle Rangitikei"
$ws.Range("T79").Value = "Rang_2d"
$ws.Range("U79").Value = "mg/L"
$ws.Range("A80").Value = "Moawhango at Waiouru"
$ws.Range("B80").Value = "DRP (Median)"
$ws.Range("C80").Value = "B"
$ws.Range("D80").Value = "2019 - 2023"
$ws.Range("E80").Value = "RepSite"
$ws.Range("F80").Value = 0.008999999999999999
$ws.Range("G80").Value = 0.0125192307692308
$ws.Range("H80").Value = 0.078
$ws.Range("I80").Value = 0.0259
$ws.Range("L80").Value = 0.008500000000000001
$ws.Range("M80").Value = 0.018
$ws.Range("N80").Value = 0.02368
$ws.Range("O80").Value = 1838908
$ws.Range("P80").Value = 5631468
$ws.Range("Q80").Value = "Ruapehu District"
$ws.Range("R80").Value = "Rangitīkei-Turakina"
$ws.Range("S80").Value = "Middle Rangitikei"
$ws.Range("T80").Value = "Rang_2d"
$ws.Range("U80").Value = "mg/L"
$ws.Range("A81").Value = "Moawhango at Waiouru"
$ws.Range("B81").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C81").Value = "A"
$ws.Range("D81").Value = "2019 - 2023"
$ws.Range("E81").Value = "RepSite"
$ws.Range("F81").Value = 0.00281
$ws.Range("G81").Value = 0.0058380233983024
$ws.Range("H81").Value = 0.0362045582721472
$ws.Range("I81").Value = 0.02613
$ws.Range("L81").Value = 0.0037
$ws.Range("M81").Value = 0.01107
$ws.Range("N81").Value = 0.01874
$ws.Range("O81").Value = 1838908
$ws.Range("P81").Value = 5631468
$ws.Range("Q81").Value = "Ruapehu District"
$ws.Range("R81").Value = "Rangitīkei-Turakina"
$ws.Range("S81").Value = "Middle Rangitikei"
$ws.Range("T81").Value = "Rang_2d"
$ws.Range("U81").Value = "mg NH4-N/L"
$ws.Range("A82").Value = "Moawhango at Waiouru"
$ws.Range("B82").Value = "Ammoniacal-N (Median)"
$ws.Range("C82").Value = "A"
$ws.Range("D82").Value = "2019 - 2023"
$ws.Range("E82").Value = "RepSite"
$ws.Range("F82").Value = 0.00281
$ws.Range("G82").Value = 0.0058380233983024
$ws.Range("H82").Value = 0.0362045582721472
$ws.Range("I82").Value = 0.02613
$ws.Range("L82").Value = 0.0037
$ws.Range("M82").Value = 0.01107
$ws.Range("N82").Value = 0.01874
$ws.Range("O82").Value = 1838908
$ws.Range("P82").Value = 5631468
$ws.Range("Q82").Value = "Ruapehu District"
$ws.Range("R82").Value = "Rangitīkei-Turakina"
$ws.Range("S82").Value = "Middle Rangitikei"
$ws.Range("T82").Value = "Rang_2d"
$ws.Range("U82").Value = "mg NH4-N/L"
$ws.Range("A83").Value = "Moawhango at Waiouru"
$ws.Range("B83").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C83").Value = "A"
$ws.Range("D83").Value = "2019 - 2023"
$ws.Range("E83").Value = "RepSite"
$ws.Range("F83").Value = 0.00128
$ws.Range("G83").Value = 0.0015319372094318
$ws.Range("H83").Value = 0.007
$ws.Range("I83").Value = 0.003
$ws.Range("L83").Value = 0.00109
$ws.Range("M83").Value = 0.002
$ws.Range("N83").Value = 0.003
$ws.Range("O83").Value = 1838908
$ws.Range("P83").Value = 5631468
$ws.Range("Q83").Value = "Ruapehu District"
$ws.Range("R83").Value = "Rangitīkei-Turakina"
$ws.Range("S83").Value = "Middle Rangitikei"
$ws.Range("T83").Value = "Rang_2d"
$ws.Range("U83").Value = "mg NO3-N/L"
$ws.Range("A84").Value = "Moawhango at Waiouru"
$ws.Range("B84").Value = "Nitrate-N (Median)"
$ws.Range("C84").Value = "A"
$ws.Range("D84").Value = "2019 - 2023"
$ws.Range("E84").Value = "RepSite"
$ws.Range("F84").Value = 0.00128
$ws.Range("G84").Value = 0.0015319372094318
$ws.Range("H84").Value = 0.007
$ws.Range("I84").Value = 0.003
$ws.Range("L84").Value = 0.00109
$ws.Range("M84").Value = 0.002
$ws.Range("N84").Value = 0.003
$ws.Range("O84").Value = 1838908
$ws.Range("P84").Value = 5631468
$ws.Range("Q84").Value = "Ruapehu District"
$ws.Range("R84").Value = "Rangitīkei-Turakina"
$ws.Range("S84").Value = "Middle Rangitikei"
$ws.Range("T84").Value = "Rang_2d"
$ws.Range("U84").Value = "mg NO3-N/L"
$ws.Range("A85").Value = "Moawhango at Waiouru"
$ws.Range("B85").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D85").Value = "2019 - 2023"
$ws.Range("E85").Value = "RepSite"
$ws.Range("F85").Value = 0.00675
$ws.Range("G85").Value = 0.0082469282836771
$ws.Range("H85").Value = 0.037
$ws.Range("I85").Value = 0.02022
$ws.Range("L85").Value = 0.007
$ws.Range("M85").Value = 0.01278
$ws.Range("N85").Value = 0.01833
$ws.Range("O85").Value = 1838908
$ws.Range("P85").Value = 5631468
$ws.Range("Q85").Value = "Ruapehu District"
$ws.Range("R85").Value = "Rangitīkei-Turakina"
$ws.Range("S85").Value = "Middle Rangitikei"
$ws.Range("T85").Value = "Rang_2d"
$ws.Range("U85").Value = "g/m3"
$ws.Range("A86").Value = "Moawhango at Waiouru"
$ws.Range("B86").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D86").Value = "2019 - 2023"
$ws.Range("E86").Value = "RepSite"
$ws.Range("F86").Value = 0.00675
$ws.Range("G86").Value = 0.0082469282836771
$ws.Range("H86").Value = 0.037
$ws.Range("I86").Value = 0.02022
$ws.Range("L86").Value = 0.007
$ws.Range("M86").Value = 0.01278
$ws.Range("N86").Value = 0.01833
$ws.Range("O86").Value = 1838908
$ws.Range("P86").Value = 5631468
$ws.Range("Q86").Value = "Ruapehu District"
$ws.Range("R86").Value = "Rangitīkei-Turakina"
$ws.Range("S86").Value = "Middle Rangitikei"
$ws.Range("T86").Value = "Rang_2d"
$ws.Range("U86").Value = "g/m3"
$ws.Range("A87").Value = "Moawhango at Waiouru"
$ws.Range("B87").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D87").Value = "2019 - 2023"
$ws.Range("E87").Value = "RepSite"
$ws.Range("F87").Value = 0.05
$ws.Range("G87").Value = 0.05575
$ws.Range("H87").Value = 0.13
$ws.Range("I87").Value = 0.115
$ws.Range("L87").Value = 0.05
$ws.Range("M87").Value = 0.07000000000000001
$ws.Range("N87").Value = 0.08599999999999999
$ws.Range("O87").Value = 1838908
$ws.Range("P87").Value = 5631468
$ws.Range("Q87").Value = "Ruapehu District"
$ws.Range("R87").Value = "Rangitīkei-Turakina"
$ws.Range("S87").Value = "Middle Rangitikei"
$ws.Range("T87").Value = "Rang_2d"
$ws.Range("U87").Value = "g/m3"
$ws.Range("A88").Value = "Moawhango at Waiouru"
$ws.Range("B88").Value = "Total Nitrogen (Median)"
$ws.Range("D88").Value = "2019 - 2023"
$ws.Range("E88").Value = "RepSite"
$ws.Range("F88").Value = 0.05
$ws.Range("G88").Value = 0.05575
$ws.Range("H88").Value = 0.13
$ws.Range("I88").Value = 0.115
$ws.Range("L88").Value = 0.05
$ws.Range("M88").Value = 0.07000000000000001
$ws.Range("N88").Value = 0.08599999999999999
$ws.Range("O88").Value = 1838908
$ws.Range("P88").Value = 5631468
$ws.Range("Q88").Value = "Ruapehu District"
$ws.Range("R88").Value = "Rangitīkei-Turakina"
$ws.Range("S88").Value = "Middle Rangitikei"
$ws.Range("T88").Value = "Rang_2d"
$ws.Range("U88").Value = "g/m3"
$ws.Range("A89").Value = "Moawhango at Waiouru"
$ws.Range("B89").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D89").Value = "2019 - 2023"
$ws.Range("E89").Value = "RepSite"
$ws.Range("F89").Value = 0.0105
$ws.Range("G89").Value = 0.0138
$ws.Range("H89").Value = 0.029
$ws.Range("I89").Value = 0.026
$ws.Range("L89").Value = 0.0095
$ws.Range("M89").Value = 0.022
$ws.Range("N89").Value = 0.0253
$ws.Range("O89").Value = 1838908
$ws.Range("P89").Value = 5631468
$ws.Range("Q89").Value = "Ruapehu District"
$ws.Range("R89").Value = "Rangitīkei-Turakina"
$ws.Range("S89").Value = "Middle Rangitikei"
$ws.Range("T89").Value = "Rang_2d"
$ws.Range("U89").Value = "g/m3"
$ws.Range("A90").Value = "Moawhango at Waiouru"
$ws.Range("B90").Value = "Total Phosphorus (Median)"
$ws.Range("D90").Value = "2019 - 2023"
$ws.Range("E90").Value = "RepSite"
$ws.Range("F90").Value = 0.0105
$ws.Range("G90").Value = 0.0138
$ws.Range("H90").Value = 0.029
$ws.Range("I90").Value = 0.026
$ws.Range("L90").Value = 0.0095
$ws.Range("M90").Value = 0.022
$ws.Range("N90").Value = 0.0253
$ws.Range("O90").Value = 1838908
$ws.Range("P90").Value = 5631468
$ws.Range("Q90").Value = "Ruapehu District"
$ws.Range("R90").Value = "Rangitīkei-Turakina"
$ws.Range("S90").Value = "Middle Rangitikei"
$ws.Range("T90").Value = "Rang_2d"
$ws.Range("U90").Value = "g/m3"
